# Apply the change described by the diff: add columns I (I0) and J (IF)
# to the worksheet, mirroring the formatting of column H and filling in
# the corresponding values for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold font, border, center/top alignment)
# from H1 onto the two new header cells so the new cells reuse the same
# cell style as the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns, keyed by row number.
$values = @{
    2  = @(8, 8)
    3  = @(3, 5)
    4  = @(9, 9)
    5  = @(6, 6)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(8, 8)
    9  = @(6, 7)
    10 = @(7, 7)
    11 = @(7, 7)
    12 = @(9, 9)
    13 = @(7, 7)
    14 = @(9, 9)
    15 = @(6, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(6, 6)
    21 = @(4, 5)
    22 = @(5, 6)
    23 = @(8, 8)
    24 = @(4, 5)
    25 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
